# The worksheet had an extra, unused column E between the data in
# columns A:D and the "Time" column which was incorrectly placed in
# column G. Remove the empty column E so the Time column (and its
# values) shift left into column F, tightening the used range from
# A1:G5 down to A1:F5.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("E").Delete()

# Restore the active cell selection recorded in the saved workbook.
$ws.Range("B20").Select()
